# Se agrega en contenido editorial otro link
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Append the new link to the existing "Contenido editorial." links cell (B6)
# and wrap the text so the long, comma separated list of links is readable.
$ws.Range("B6").Value = "http://inovom.pruebab2b.com/temas/+113492,     http://inovom.pruebab2b.com/temas/+103286,   http://inovom.pruebab2b.com/temas/+3084014, http://inovom.pruebab2b.com/temas/+112326"
$ws.Range("B6").WrapText = $true

# Clear out the stray "PENDIENTE" markers that were left in column C for rows
# that no longer need them.
$ws.Range("C29").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").ClearContents()

# Update the saved view/selection state of the sheet.
$ws.Activate()
$ws.Range("D23").Select()
